# Insert a new weekly record row into the Coliflor dataset.
# This pushes the existing rows 1092..1139 down to 1093..1140 and
# populates the newly freed row 1092 with the new data point.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 1092 (shifts 1092:1139 -> 1093:1140,
# and copies formatting/number-format down from the row above, matching
# Excel's default "insert row" behavior).
$ws.Rows.Item(1092).Insert()

# Fill in the constant columns (same for every record in this block).
$ws.Range("A1092").Value = 3
$ws.Range("B1092").Value = "Femacal de La Calera"
$ws.Range("C1092").Value = "Coquimbo"
$ws.Range("D1092").Value = 45147
$ws.Range("E1092").Value = 5
$ws.Range("F1092").Value = 100112008
$ws.Range("G1092").Value = "Coliflor"
$ws.Range("H1092").Value = "Sin especificar"
$ws.Range("I1092").Value = "Primera"
$ws.Range("J1092").Value = 2800
$ws.Range("K1092").Value = 700
$ws.Range("L1092").Value = 750
$ws.Range("M1092").Value = 721
$ws.Range("N1092").Value = "$/unidad"
$ws.Range("O1092").Value = "Provincia de Quillota"
$ws.Range("P1092").Value = 721
$ws.Range("Q1092").Value = 1
$ws.Range("R1092").Value = "Hortaliza"

# Keep the date column's number format consistent with the rest of column D.
$ws.Range("D1092").NumberFormat = $ws.Range("D1093").NumberFormat()
